# Apply updated market-board price/profit figures to the Leve tables.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 470.64706
$ws.Range("I33").Value = 429.83334
$ws.Range("K33").Value = 429.83334
$ws.Range("M33").Value = -200.83334

$ws.Range("H70").Value = 6655.5557
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 7237.5
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 21712.5
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -22252.5

$ws.Range("H73").Value = 6655.5557
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 7237.5
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 21712.5
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -23584.5

$ws.Range("H94").Value = 2010.7142
$ws.Range("I94").Value = 2010.7142
$ws.Range("K94").Value = 2010.7142
$ws.Range("M94").Value = -1559.7142

$ws.Range("H99").Value = 500
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H103").Value = 924.6667
$ws.Range("I103").Value = 856.125
$ws.Range("J103").Value = 1061.75
$ws.Range("K103").Value = 2568.375
$ws.Range("L103").Value = 3185.25
$ws.Range("M103").Value = -1982.375
$ws.Range("N103").Value = -4357.25

$ws.Range("H113").Value = 1468
$ws.Range("I113").Value = 1468
$ws.Range("K113").Value = 1468
$ws.Range("M113").Value = 1786

$ws.Range("H138").Value = 4249.8335
$ws.Range("J138").Value = 4600.1
$ws.Range("L138").Value = 13800.3
$ws.Range("N138").Value = -24080.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 34975
$ws.Range("J24").Value = 34975
$ws.Range("L24").Value = 34975
$ws.Range("N24").Value = -35723

$ws.Range("H32").Value = 4534.5454
$ws.Range("I32").Value = 2344.3572
$ws.Range("K32").Value = 2344.3572
$ws.Range("M32").Value = -2057.3572

$ws.Range("H45").Value = 6544.3335
$ws.Range("I45").Value = 7771.2856
$ws.Range("J45").Value = 2250
$ws.Range("K45").Value = 7771.2856
$ws.Range("L45").Value = 2250
$ws.Range("M45").Value = -7394.2856
$ws.Range("N45").Value = -3004

$ws.Range("H61").Value = 250006750
$ws.Range("J61").Value = 8998
$ws.Range("L61").Value = 8998
$ws.Range("N61").Value = -9422

$ws.Range("H88").Value = 3999.5
$ws.Range("I88").Value = 3949
$ws.Range("J88").Value = 4024.75
$ws.Range("K88").Value = 3949
$ws.Range("L88").Value = 4024.75
$ws.Range("M88").Value = -3543
$ws.Range("N88").Value = -4836.75

$ws.Range("H91").Value = 3999.5
$ws.Range("I91").Value = 3949
$ws.Range("J91").Value = 4024.75
$ws.Range("K91").Value = 3949
$ws.Range("L91").Value = 4024.75
$ws.Range("M91").Value = -2545
$ws.Range("N91").Value = -6832.75

$ws.Range("H100").Value = 34975
$ws.Range("J100").Value = 34975
$ws.Range("L100").Value = 34975
$ws.Range("N100").Value = -37139

$ws.Range("H102").Value = 2684.375
$ws.Range("I102").Value = 2650.4614
$ws.Range("K102").Value = 2650.4614
$ws.Range("M102").Value = -1028.4614

$ws.Range("H132").Value = 2858613
$ws.Range("I132").Value = 3031671.2
$ws.Range("J132").Value = 3151.5
$ws.Range("K132").Value = 9095013.600000001
$ws.Range("L132").Value = 9454.5
$ws.Range("M132").Value = -9092483.600000001
$ws.Range("N132").Value = -14514.5

$ws.Range("H136").Value = 250006750
$ws.Range("J136").Value = 8998
$ws.Range("L136").Value = 26994
$ws.Range("N136").Value = -32094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1833.9
$ws.Range("I20").Value = 1707.6
$ws.Range("J20").Value = 1960.2
$ws.Range("K20").Value = 1707.6
$ws.Range("L20").Value = 1960.2
$ws.Range("M20").Value = -1460.6
$ws.Range("N20").Value = -2454.2

$ws.Range("H94").Value = 2202.2273
$ws.Range("I94").Value = 2232.45
$ws.Range("K94").Value = 2232.45
$ws.Range("M94").Value = -1781.45

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10229
$ws.Range("I31").Value = 36500
$ws.Range("J31").Value = 1472
$ws.Range("K31").Value = 36500
$ws.Range("L31").Value = 1472
$ws.Range("M31").Value = -36205
$ws.Range("N31").Value = -2062

$ws.Range("H34").Value = 10229
$ws.Range("I34").Value = 36500
$ws.Range("J34").Value = 1472
$ws.Range("K34").Value = 36500
$ws.Range("L34").Value = 1472
$ws.Range("M34").Value = -36298
$ws.Range("N34").Value = -1876

$ws.Range("H58").Value = 22736046
$ws.Range("I58").Value = 26325150
$ws.Range("J58").Value = 5061.6665
$ws.Range("K58").Value = 26325150
$ws.Range("L58").Value = 5061.6665
$ws.Range("M58").Value = -26324947
$ws.Range("N58").Value = -5467.6665

$ws.Range("H132").Value = 30305142
$ws.Range("I132").Value = 32260216
$ws.Range("K132").Value = 96780648
$ws.Range("M132").Value = -96778118

$ws.Range("H134").Value = 27780512
$ws.Range("I134").Value = 35716908
$ws.Range("K134").Value = 107150724
$ws.Range("M134").Value = -107148189

$ws.Range("H136").Value = 22736046
$ws.Range("I136").Value = 26325150
$ws.Range("J136").Value = 5061.6665
$ws.Range("K136").Value = 78975450
$ws.Range("L136").Value = 15184.9995
$ws.Range("M136").Value = -78972900
$ws.Range("N136").Value = -20284.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3431.8572
$ws.Range("I115").Value = 3341
$ws.Range("J115").Value = 3500
$ws.Range("K115").Value = 10023
$ws.Range("L115").Value = 10500
$ws.Range("M115").Value = -8848
$ws.Range("N115").Value = -12850

$ws.Range("H122").Value = 727.1539
$ws.Range("I122").Value = 698.4
$ws.Range("K122").Value = 6285.599999999999
$ws.Range("M122").Value = -3835.599999999999

$ws.Range("H137").Value = 2180.5715
$ws.Range("I137").Value = 1952.8
$ws.Range("K137").Value = 5858.4
$ws.Range("M137").Value = -758.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4027.0908
$ws.Range("I70").Value = 3057.1428
$ws.Range("J70").Value = 5724.5
$ws.Range("K70").Value = 3057.1428
$ws.Range("L70").Value = 5724.5
$ws.Range("M70").Value = -2787.1428
$ws.Range("N70").Value = -6264.5

$ws.Range("H73").Value = 4027.0908
$ws.Range("I73").Value = 3057.1428
$ws.Range("J73").Value = 5724.5
$ws.Range("K73").Value = 3057.1428
$ws.Range("L73").Value = 5724.5
$ws.Range("M73").Value = -2121.1428
$ws.Range("N73").Value = -7596.5

$ws.Range("H102").Value = 2129.3684
$ws.Range("I102").Value = 2129.3684
$ws.Range("K102").Value = 2129.3684
$ws.Range("M102").Value = -507.3683999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3012.5
$ws.Range("I22").Value = 3329.6
$ws.Range("J22").Value = 2484
$ws.Range("K22").Value = 3329.6
$ws.Range("L22").Value = 2484
$ws.Range("M22").Value = -3034.6
$ws.Range("N22").Value = -3074

$ws.Range("H27").Value = 3012.5
$ws.Range("I27").Value = 3329.6
$ws.Range("J27").Value = 2484
$ws.Range("K27").Value = 3329.6
$ws.Range("L27").Value = 2484
$ws.Range("M27").Value = -3222.6
$ws.Range("N27").Value = -2698

$ws.Range("H46").Value = 768
$ws.Range("I46").Value = 755.2857
$ws.Range("K46").Value = 755.2857
$ws.Range("M46").Value = -567.2857

$ws.Range("H55").Value = 502.31033
$ws.Range("I55").Value = 342.05554
$ws.Range("J55").Value = 764.5454999999999
$ws.Range("K55").Value = 342.05554
$ws.Range("L55").Value = 764.5454999999999
$ws.Range("M55").Value = -169.05554
$ws.Range("N55").Value = -1110.5455

$ws.Range("H93").Value = 2132.2307
$ws.Range("I93").Value = 1861.9
$ws.Range("J93").Value = 3033.3333
$ws.Range("K93").Value = 1861.9
$ws.Range("L93").Value = 3033.3333
$ws.Range("M93").Value = -613.9000000000001
$ws.Range("N93").Value = -5529.3333

$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1745
$ws.Range("I96").Value = 1660
$ws.Range("K96").Value = 1660
$ws.Range("M96").Value = -287

$ws.Range("H116").Value = 66220
$ws.Range("J116").Value = 66220
$ws.Range("L116").Value = 66220
$ws.Range("N116").Value = -75398

$ws.Range("H126").Value = 1934.32
$ws.Range("I126").Value = 1799.8948
$ws.Range("J126").Value = 2360
$ws.Range("K126").Value = 5399.6844
$ws.Range("L126").Value = 7080
$ws.Range("M126").Value = -2929.6844
$ws.Range("N126").Value = -12020

$ws.Range("H136").Value = 23810692
$ws.Range("I136").Value = 27778806
$ws.Range("J136").Value = 1999.6666
$ws.Range("K136").Value = 83336418
$ws.Range("L136").Value = 5998.9998
$ws.Range("M136").Value = -83333868
$ws.Range("N136").Value = -11098.9998
